$wb = $excel.ActiveWorkbook

$configWs = $wb.Worksheets.Item("config")
$configWs.Activate()

$configWs.Range("D1").Value = "commodity"
$configWs.Range("E1").Value = "year"
$configWs.Range("D2").Value = "light"
$configWs.Range("E2").Value = 700
$configWs.Range("E3").Value = 710
$configWs.Range("E4").Value = 720

$configWs.Range("E7").Select()
